$d = $word.ActiveDocument

# Locate the unique run containing the InsuranceAmount B/C / Strata check.
$needle = [string]::Format("{0}tr if InsuranceAmount != {1}B/C{2} and InsuranceAmount != {1}Strata{2} %{3}", "{%", [char]0x201C, [char]0x201D, "}")

$full = $d.Content
$found = $full.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text"
}

$segStart = $full.Start

# Compute the start/end offsets of the word "Strata" within the found range.
$strataLocal = $needle.IndexOf("Strata")
$strataStart = $segStart + $strataLocal
$strataEnd = $strataStart + 6

$target = $d.Range($strataStart, $strataEnd)

# Force Word to split this span into its own run (distinct from its
# neighbours) by toggling a character-formatting property across the
# text replacement, then reverting that property on the now-isolated run.
$target.Bold = $true
$target.Text = "STRATA"
$target.Bold = $false
